$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("sigma_010")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 27.76911232435696
$ws.Range("C2").Value = 30.95469958623617
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 27.77218200925297
$ws.Range("C3").Value = 30.94594025855241
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 27.77961035142656
$ws.Range("C4").Value = 30.95196471224806
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 27.77726214367717
$ws.Range("C5").Value = 30.9284795414577
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 27.77727439371471
$ws.Range("C6").Value = 30.94385957015226
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 27.78616386283362
$ws.Range("C7").Value = 30.9421712131131
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 27.75530679326049
$ws.Range("C8").Value = 30.92070938909953
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 27.75427856400638
$ws.Range("C9").Value = 30.91092999684085
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 27.7901244056256
$ws.Range("C10").Value = 30.9329328084545
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 27.74684977610189
$ws.Range("C11").Value = 30.91411999277354
$ws.Range("B12").Value = 27.77081646242563
$ws.Range("C12").Value = 30.93458070689281

$ws = $wb.Worksheets.Item("sigma_025")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 19.65532595139404
$ws.Range("C2").Value = 27.72216175425647
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 19.67921997613642
$ws.Range("C3").Value = 27.7557374228078
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 19.66943412625176
$ws.Range("C4").Value = 27.72065624395826
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 19.66423290171531
$ws.Range("C5").Value = 27.74387410776457
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 19.66589164591669
$ws.Range("C6").Value = 27.76003178999528
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 19.66379350465722
$ws.Range("C7").Value = 27.73792074456587
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 19.67290162511125
$ws.Range("C8").Value = 27.75306045787445
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 19.68033866088162
$ws.Range("C9").Value = 27.74286629429223
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 19.66450176155188
$ws.Range("C10").Value = 27.77910290379109
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 19.67037031089562
$ws.Range("C11").Value = 27.78218657562942
$ws.Range("B12").Value = 19.66860104645118
$ws.Range("C12").Value = 27.74975982949355

$ws = $wb.Worksheets.Item("sigma_050")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 14.57274643124718
$ws.Range("C2").Value = 22.39794536257531
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 14.5615505613959
$ws.Range("C3").Value = 22.41449991866852
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 14.57925255383172
$ws.Range("C4").Value = 22.43477074005446
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 14.58060040727101
$ws.Range("C5").Value = 22.45316896784264
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 14.56082076980098
$ws.Range("C6").Value = 22.44469786905776
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 14.56324355222458
$ws.Range("C7").Value = 22.43152196630176
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 14.54821185866155
$ws.Range("C8").Value = 22.45531765322306
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 14.57407853556555
$ws.Range("C9").Value = 22.45838247355496
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 14.58730029836311
$ws.Range("C10").Value = 22.40950463928297
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 14.56180248557641
$ws.Range("C11").Value = 22.46903412550548
$ws.Range("B12").Value = 14.5689607453938
$ws.Range("C12").Value = 22.43688437160669
